{"js": "// Update project aim paragraph text: \"\u0422\u0435\u0441\u0442\u043e\u0432\u0430 \u0446\u0456\u043b\u044c\" -> \"\u0426\u0435 \u043f\u0440\u043e\u0435\u043a\u0442 \u0434\u043b\u044f \u0442\u0435\u0441\u0442\u0443\"\nconst aimResults = context.document.body.search(\"\u0422\u0435\u0441\u0442\u043e\u0432\u0430 \u0446\u0456\u043b\u044c\", { matchCase: true });\naimResults.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < aimResults.items.length; i++) {\n  aimResults.items[i].insertText(\"\u0426\u0435 \u043f\u0440\u043e\u0435\u043a\u0442 \u0434\u043b\u044f \u0442\u0435\u0441\u0442\u0443\", \"Replace\");\n}\nawait context.sync();\n\n// Update glossary table: \"\u041e\u0437\u043d\u0430\u0447\u0435\u043d\u043d\u044f 1\" -> \"\u0417\u043d\u0430\u0447\u0435\u043d\u043d\u044f 1\" and append two new rows.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst glossaryTable = tables.items[0];\nconst rows = glossaryTable.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Fix the existing definition cell text in the first data row (row index 1,\n// since row index 0 is the header \"Term\"/\"Definition\").\nconst defCell = glossaryTable.getCell(1, 1);\ndefCell.value = \"\u0417\u043d\u0430\u0447\u0435\u043d\u043d\u044f 1\";\n\n// Append two new rows with the same shape as the existing data row.\nglossaryTable.addRows(\"End\", 2, [\n  [\"\u0422\u0435\u0440\u043c\u0456\u043d 2\", \"\u0417\u043d\u0430\u0447\u0435\u043d\u043d\u044f 2\"],\n  [\"\u0422\u0435\u0440\u043c\u0456\u043d 3\", \"\u0417\u043d\u0430\u0447\u0435\u043d\u043d\u044f 3\"]\n]);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the \"Project aim\" paragraph text.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"\u0422\u0435\u0441\u0442\u043e\u0432\u0430 \u0446\u0456\u043b\u044c\"\n$find.Replacement.Text = \"\u0426\u0435 \u043f\u0440\u043e\u0435\u043a\u0442 \u0434\u043b\u044f \u0442\u0435\u0441\u0442\u0443\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# Update the Glossary table: fix \"\u041e\u0437\u043d\u0430\u0447\u0435\u043d\u043d\u044f 1\" -> \"\u0417\u043d\u0430\u0447\u0435\u043d\u043d\u044f 1\" and append two\n# more Term/Definition rows with the same formatting as the existing row.\n$tbl = $d.Tables.Item(1)\n$tbl.Cell(2, 2).Range.Text = \"\u0417\u043d\u0430\u0447\u0435\u043d\u043d\u044f 1\"\n\n$row2 = $tbl.Rows.Add()\n$row2.Cells.Item(1).Range.Text = \"\u0422\u0435\u0440\u043c\u0456\u043d 2\"\n$row2.Cells.Item(2).Range.Text = \"\u0417\u043d\u0430\u0447\u0435\u043d\u043d\u044f 2\"\n\n$row3 = $tbl.Rows.Add()\n$row3.Cells.Item(1).Range.Text = \"\u0422\u0435\u0440\u043c\u0456\u043d 3\"\n$row3.Cells.Item(2).Range.Text = \"\u0417\u043d\u0430\u0447\u0435\u043d\u043d\u044f 3\"\n"}
